$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Ebola GP1 (Zaire)").Name = "Zaire Ebola GP"
$wb.Worksheets.Item("Ebola GP1 (Sudan)").Name = "Sudan Ebola GP"
$wb.Worksheets.Item("Ebola NP (Zaire)").Name = "Zaire Ebola NP"
$wb.Worksheets.Item("Ebola NP (Sudan)").Name = "Sudan Ebola NP"
$wb.Worksheets.Item("Burkholderia HCP1").Name = "Burkholderia Hcp1"
